$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: "How did you feel before testing the game?" / "How did
# you feel immediately afterwards?" (two separate paragraphs) become
# a single paragraph "How did you like the game?", split across three
# runs ("How did you " / "like the game" / "?") that share identical
# run formatting (sz=24, szCs=24).
# ------------------------------------------------------------------

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13,[char]7) -eq "How did you feel before testing the game?") {
        $target = $p
        break
    }
}

$pStart = $target.Range.Start
$midStart = $pStart + 12          # end of "How did you "
$midEnd = $midStart + 28          # end of "feel before testing the game" (before "?")

# Replace the middle portion ("feel before testing the game") with
# the new wording ("like the game").
$midRng = $d.Range($midStart, $midEnd)
$midRng.Text = "like the game"

# Force the three segments to stay as distinct runs (re-stamping the
# font size to its existing value triggers a run split without
# changing any visible formatting).
$r1 = $d.Range($pStart, $midStart)
$r2 = $d.Range($midStart, $midStart + 13)
$r3 = $d.Range($midStart + 13, $midStart + 14)
$r1.Font.Size = 11
$r1.Font.Size = 12
$r2.Font.Size = 11
$r2.Font.Size = 12
$r3.Font.Size = 11
$r3.Font.Size = 12

# Remove the now-redundant "How did you feel immediately afterwards?"
# paragraph that used to follow it.
$nextPara = $target.Next()
$nextPara.Range.Delete()

# ------------------------------------------------------------------
# Change 2: add three new survey questions right after "Were there
# any times where the difficulty felt unfair?".
# ------------------------------------------------------------------

$anchor = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13,[char]7) -eq "Were there any times where the difficulty felt unfair?") {
        $anchor = $p
        break
    }
}

$newQuestions = @(
    "Was the control scheme intuitive?",
    "Was the heads-up display intuitive?",
    "Was the combat viscerally satisfying?"
)

foreach ($q in $newQuestions) {
    $anchor.Range.InsertParagraphAfter()
    $anchor = $anchor.Next()
    $anchor.Range.Text = $q
}
